$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Sicil) values for rows 3 through 99 to continue the
# date-serial sequence that already starts at C2 (44142).
for ($row = 3; $row -le 99; $row++) {
    $ws.Cells.Item($row, 3).Value = 44142 + ($row - 2)
}

# Update the active selection to E6
$ws.Range("E6").Select()
